$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1910.8125
$ws.Range("I18").Value = 1840.9286
$ws.Range("J18").Value = 2400
$ws.Range("K18").Value = 1840.9286
$ws.Range("L18").Value = 2400
$ws.Range("M18").Value = -1556.9286
$ws.Range("N18").Value = -2968
$ws.Range("H51").Value = 6948.6
$ws.Range("I51").Value = 6833.1665
$ws.Range("K51").Value = 6833.1665
$ws.Range("M51").Value = -6349.1665
$ws.Range("H62").Value = 5484.25
$ws.Range("I62").Value = 3018.7
$ws.Range("K62").Value = 3018.7
$ws.Range("M62").Value = -2394.7
$ws.Range("H65").Value = 5484.25
$ws.Range("I65").Value = 3018.7
$ws.Range("K65").Value = 15093.5
$ws.Range("M65").Value = -11973.5
$ws.Range("H100").Value = 2429
$ws.Range("I100").Value = 2143.889
$ws.Range("K100").Value = 2143.889
$ws.Range("M100").Value = -1602.889
$ws.Range("H104").Value = 207.4
$ws.Range("I104").Value = 217.11111
$ws.Range("K104").Value = 651.3333299999999
$ws.Range("M104").Value = 1095.66667
$ws.Range("H111").Value = 11115813
$ws.Range("I111").Value = 22230656
$ws.Range("J111").Value = 970
$ws.Range("K111").Value = 66691968
$ws.Range("L111").Value = 2910
$ws.Range("M111").Value = -66688901
$ws.Range("N111").Value = -9044
$ws.Range("H116").Value = 5061.2104
$ws.Range("J116").Value = 5755.3076
$ws.Range("L116").Value = 5755.3076
$ws.Range("N116").Value = -12639.3076
$ws.Range("H131").Value = 4197.609
$ws.Range("I131").Value = 2269.7334
$ws.Range("K131").Value = 6809.2002
$ws.Range("M131").Value = -1769.2002
$ws.Range("H137").Value = 95980.48
$ws.Range("J137").Value = 3748.5
$ws.Range("L137").Value = 11245.5
$ws.Range("N137").Value = -16345.5
$ws.Range("H138").Value = 5457.0527
$ws.Range("J138").Value = 5542.9375
$ws.Range("L138").Value = 16628.8125
$ws.Range("N138").Value = -26908.8125
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1985690.2
$ws.Range("I110").Value = 5558037.5
$ws.Range("K110").Value = 5558037.5
$ws.Range("M110").Value = -5555992.5
$ws.Range("H122").Value = 536272.7
$ws.Range("I122").Value = 1614.5
$ws.Range("K122").Value = 4843.5
$ws.Range("M122").Value = -2393.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2020951.2
$ws.Range("I94").Value = 2755376.8
$ws.Range("J94").Value = 1281.1666
$ws.Range("K94").Value = 2755376.8
$ws.Range("L94").Value = 1281.1666
$ws.Range("M94").Value = -2754925.8
$ws.Range("N94").Value = -2183.1666
$ws.Range("H132").Value = 85000
$ws.Range("J132").Value = 85000
$ws.Range("L132").Value = 85000
$ws.Range("N132").Value = -95120
$ws.Range("H134").Value = 11650.963
$ws.Range("I134").Value = 11253.728
$ws.Range("K134").Value = 33761.18399999999
$ws.Range("M134").Value = -31226.18399999999
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H140").Value = 88995
$ws.Range("J140").Value = 88995
$ws.Range("L140").Value = 88995
$ws.Range("N140").Value = -99355

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 29998.334
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 29998.334
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 29998.334
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -31470.334
$ws.Range("H61").Value = 29998.334
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 29998.334
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 29998.334
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -30694.334
$ws.Range("H62").Value = 2624.4375
$ws.Range("I62").Value = 2014.7693
$ws.Range("K62").Value = 2014.7693
$ws.Range("M62").Value = -1390.7693
$ws.Range("H65").Value = 2624.4375
$ws.Range("I65").Value = 2014.7693
$ws.Range("K65").Value = 10073.8465
$ws.Range("M65").Value = -6953.8465
$ws.Range("H140").Value = 84999
$ws.Range("J140").Value = 84999
$ws.Range("L140").Value = 84999
$ws.Range("N140").Value = -95359

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26442794
$ws.Range("I4").Value = 10845039
$ws.Range("J4").Value = 37584050
$ws.Range("K4").Value = 32535117
$ws.Range("L4").Value = 112752150
$ws.Range("M4").Value = -32535005
$ws.Range("N4").Value = -112752374
$ws.Range("H59").Value = 2882.8
$ws.Range("I59").Value = 2905
$ws.Range("J59").Value = 2877.25
$ws.Range("K59").Value = 8715
$ws.Range("L59").Value = 8631.75
$ws.Range("M59").Value = -8175
$ws.Range("N59").Value = -9711.75
$ws.Range("H60").Value = 2994
$ws.Range("J60").Value = 2994
$ws.Range("L60").Value = 8982
$ws.Range("N60").Value = -9484
$ws.Range("H61").Value = 168.5
$ws.Range("I61").Value = 124.666664
$ws.Range("K61").Value = 373.999992
$ws.Range("M61").Value = -158.999992
$ws.Range("H107").Value = 1313.0714
$ws.Range("I107").Value = 1002.5
$ws.Range("K107").Value = 3007.5
$ws.Range("M107").Value = -1087.5
$ws.Range("H140").Value = 2488.75
$ws.Range("I140").Value = 2488.75
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 7466.25
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2286.25
$ws.Range("N140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25005688
$ws.Range("I70").Value = 50004124
$ws.Range("J70").Value = 7249.75
$ws.Range("K70").Value = 50004124
$ws.Range("L70").Value = 7249.75
$ws.Range("M70").Value = -50003854
$ws.Range("N70").Value = -7789.75
$ws.Range("H73").Value = 25005688
$ws.Range("I73").Value = 50004124
$ws.Range("J73").Value = 7249.75
$ws.Range("K73").Value = 50004124
$ws.Range("L73").Value = 7249.75
$ws.Range("M73").Value = -50003188
$ws.Range("N73").Value = -9121.75
$ws.Range("H97").Value = 1036516.25
$ws.Range("J97").Value = 1144.5
$ws.Range("L97").Value = 1144.5
$ws.Range("N97").Value = -2136.5
$ws.Range("H107").Value = 1236.909
$ws.Range("J107").Value = 834.3333
$ws.Range("L107").Value = 834.3333
$ws.Range("N107").Value = -4674.3333
$ws.Range("H122").Value = 991490.4399999999
$ws.Range("I122").Value = 1274064.4
$ws.Range("J122").Value = 2481.5
$ws.Range("K122").Value = 3822193.2
$ws.Range("L122").Value = 7444.5
$ws.Range("M122").Value = -3819743.2
$ws.Range("N122").Value = -12344.5
$ws.Range("H126").Value = 3682524.2
$ws.Range("I126").Value = 1820656.6
$ws.Range("K126").Value = 5461969.800000001
$ws.Range("M126").Value = -5459499.800000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 29500
$ws.Range("J136").Value = 29000
$ws.Range("L136").Value = 87000
$ws.Range("N136").Value = -92100

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8028.2856
$ws.Range("I7").Value = 4925
$ws.Range("J7").Value = 12166
$ws.Range("K7").Value = 4925
$ws.Range("L7").Value = 12166
$ws.Range("M7").Value = -4813
$ws.Range("N7").Value = -12390
$ws.Range("H46").Value = 7281.8823
$ws.Range("J46").Value = 11999.6
$ws.Range("L46").Value = 11999.6
$ws.Range("N46").Value = -12375.6
$ws.Range("H122").Value = 5866.615
$ws.Range("I122").Value = 3012.7144
$ws.Range("J122").Value = 9196.166999999999
$ws.Range("K122").Value = 9038.143199999999
$ws.Range("L122").Value = 27588.501
$ws.Range("M122").Value = -6588.143199999999
$ws.Range("N122").Value = -32488.501
$ws.Range("H126").Value = 8028.2856
$ws.Range("I126").Value = 4925
$ws.Range("J126").Value = 12166
$ws.Range("K126").Value = 14775
$ws.Range("L126").Value = 36498
$ws.Range("M126").Value = -12305
$ws.Range("N126").Value = -41438
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2540.2144
$ws.Range("I126").Value = 2380.4167
$ws.Range("K126").Value = 7141.250100000001
$ws.Range("M126").Value = -4671.250100000001
$ws.Range("H132").Value = 26612146
$ws.Range("I132").Value = 50012536
$ws.Range("J132").Value = 611711.4399999999
$ws.Range("K132").Value = 150037608
$ws.Range("L132").Value = 1835134.32
$ws.Range("M132").Value = -150035078
$ws.Range("N132").Value = -1840194.32
